# The deck's design ("Integral") theme colours are swapped back to the
# stock "Office Theme" colour values (dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink). This mirrors the commit, which exchanges the content of
# ppt/theme/theme1.xml ("Office Theme") and ppt/theme/theme2.xml
# ("Integral") — the slide master (and therefore every slide/layout)
# keeps using the same theme part, but that part's 12 scheme colours
# need to become the "Office Theme" values.

$p = $ppt.ActivePresentation

function HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target values: the "Office Theme" colour scheme (previously living in
# ppt/theme/theme1.xml) that should now be applied to the presentation's
# active theme (ppt/theme/theme2.xml), in clrScheme slot order.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $tcs.Colors($i).RGB = HexToRGB $officeThemeColors[$i - 1]
}
